$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): per-question marks
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 ("Total"): totals derived from the marking scheme
$ws.Range("B12").Value = 105
$ws.Range("C12").Value = -2.4
$ws.Range("E12").Value = "102.6/140"
